$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1, 0, 1, 1, 1, 0, 3, 0, 1, 2, 0, 0, 0, 1, 0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
